$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "link_lattes"
$ws.Range("B1").Value = "PPG"
$ws.Range("C1").Value = "nome"
$ws.Range("D1").Value = "ultima_atualizacao"
$ws.Range("E1").Value = "endereco_prof"
$ws.Range("F1").Value = "ano_ultima_formacao"
$ws.Range("G1").Value = "formacao_titulo"
$ws.Range("H1").Value = "formacao_ies"
$ws.Range("I1").Value = "ultimo_vinculo_ies"
$ws.Range("J1").Value = "prod_artigos_completos"
